$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: Coin (B), Link (C), Price (D), Volume(1h) (E)
# Price values are forced to Text via NumberFormat "@" so numeric-looking
# strings like "1.00" are not coerced to numbers by Excel, then the style
# is reset back to Normal so no stray formatting is left behind.
$updates = @(
    @{Row=2; D='42.444.19'; E='  -2.56%  '},
    @{Row=3; D='2.222.92'; E='  -2.11%  '},
    @{Row=4; D='1.00'; E='  +0.26%  '},
    @{Row=5; D='109.64'; E='  -8.04%  '},
    @{Row=6; D='296.60'; E='  +11.56%  '},
    @{Row=7; E='  -3.54%  '},
    @{Row=8; E='  -0.20%  '},
    @{Row=9; D='0.604'; E='  -2.68%  '},
    @{Row=10; D='43.99'; E='  -7.73%  '},
    @{Row=11; E='  -3.34%  '},
    @{Row=12; D='54.30'; E='  +0.15%  '},
    @{Row=13; D='8.78'; E='  -4.72%  '},
    @{Row=14; E='  +10.77%  '},
    @{Row=15; E='  -2.55%  '},
    @{Row=16; E='  -2.41%  '},
    @{Row=17; D='2.555.76'; E='  -2.16%  '},
    @{Row=18; D='2.241.39'; E='  -1.05%  '},
    @{Row=19; D='42.314.47'; E='  -2.87%  '},
    @{Row=20; D='7.37'; E='  +7.30%  '},
    @{Row=21; E='  -3.96%  '},
    @{Row=22; D='72.27'; E='  +0.19%  '},
    @{Row=23; D='3.49'; E='  +21.20%  '},
    @{Row=24; D='2.32'; E='  -3.66%  '},
    @{Row=25; D='228.80'; E='  -3.01%  '},
    @{Row=26; E='  -3.79%  '},
    @{Row=27; D='11.70'; E='  -2.89%  '},
    @{Row=28; E='  -1.70%  '},
    @{Row=29; E='  -0.73%  '},
    @{Row=30; D='38.22'; E='  -8.56%  '},
    @{Row=31; E='  -5.45%  '},
    @{Row=32; D='173.76'; E='  +0.94%  '},
    @{Row=33; E='  -2.73%  '},
    @{Row=34; D='0.0899'; E='  -1.94%  '},
    @{Row=35; D='5.69'; E='  -0.68%  '},
    @{Row=36; E='  +11.82%  '},
    @{Row=37; D='4.37'; E='  +4.16%  '},
    @{Row=38; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.125'; E='  -3.89%  '},
    @{Row=39; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.0378'; E='  -1.51%  '},
    @{Row=41; E='  -5.38%  '},
    @{Row=42; D='71.88'; E='  -2.71%  '},
    @{Row=43; D='0.235'; E='  -1.19%  '},
    @{Row=44; E='  -0.08%  '},
    @{Row=45; E='  -8.13%  '},
    @{Row=46; E='  -4.47%  '},
    @{Row=48; D='1.32'; E='  +3.84%  '},
    @{Row=49; D='103.65'; E='  +1.68%  '},
    @{Row=50; E='  -1.33%  '},
    @{Row=51; E='  +5.30%  '}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        $dCell = $ws.Cells.Item($r, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
